# Auto-generated edit script to update cryptos worksheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '37.575.07'
Set-TextValue 'E2' '  +2.26%  '
Set-TextValue 'D3' '2.044.15'
Set-TextValue 'E3' '  +3.46%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '255.44'
Set-TextValue 'E5' '  +4.22%  '
Set-TextValue 'D6' '0.623'
Set-TextValue 'E6' '  -1.06%  '
Set-TextValue 'E7' '  -5.12%  '
Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.390'
Set-TextValue 'E9' '  +2.28%  '
Set-TextValue 'D10' '57.16'
Set-TextValue 'E10' '  -0.20%  '
Set-TextValue 'D11' '0.0803'
Set-TextValue 'E11' '  +0.52%  '
Set-TextValue 'D12' '0.103'
Set-TextValue 'E12' '  -0.31%  '
Set-TextValue 'D13' '14.90'
Set-TextValue 'E13' '  +3.51%  '
Set-TextValue 'D14' '2.344.75'
Set-TextValue 'E14' '  +3.41%  '
Set-TextValue 'E15' '  -2.23%  '
Set-TextValue 'D16' '21.51'
Set-TextValue 'E16' '  -2.26%  '
Set-TextValue 'D17' '5.41'
Set-TextValue 'E17' '  +0.05%  '
Set-TextValue 'D18' '2.046.37'
Set-TextValue 'E18' '  +3.48%  '
Set-TextValue 'D19' '37.495.02'
Set-TextValue 'E19' '  +2.21%  '
Set-TextValue 'D20' '70.24'
Set-TextValue 'E20' '  +0.20%  '
Set-TextValue 'D21' '0.0₃0860'
Set-TextValue 'E21' '  +0.08%  '
Set-TextValue 'D22' '5.27'
Set-TextValue 'E22' '  +2.76%  '
Set-TextValue 'D23' '229.73'
Set-TextValue 'E23' '  -0.16%  '
Set-TextValue 'D24' '2.67'
Set-TextValue 'E24' '  +7.84%  '
Set-TextValue 'E25' '  -0.05%  '
Set-TextValue 'E26' '  -0.95%  '
Set-TextValue 'D27' '0.141'
Set-TextValue 'E27' '  -3.28%  '
Set-TextValue 'D28' '9.19'
Set-TextValue 'E28' '  -0.48%  '
Set-TextValue 'D29' '163.88'
Set-TextValue 'E29' '  +0.68%  '
Set-TextValue 'D30' '20.01'
Set-TextValue 'E30' '  +2.81%  '
Set-TextValue 'E31' '  +0.11%  '
Set-TextValue 'E32' '  -0.16%  '
Set-TextValue 'D33' '4.79'
Set-TextValue 'E33' '  -0.72%  '
Set-TextValue 'D34' '0.0668'
Set-TextValue 'E34' '  +7.51%  '
Set-TextValue 'E35' '  +0.48%  '
Set-TextValue 'E36' '  +9.69%  '
Set-TextValue 'D37' '3.50'
Set-TextValue 'E37' '  +4.63%  '
Set-TextValue 'E38' '  +0.05%  '
Set-TextValue 'E39' '  +2.22%  '
Set-TextValue 'E40' '  -0.48%  '
Set-TextValue 'E41' '  +4.16%  '
Set-TextValue 'D42' '0.0973'
Set-TextValue 'E42' '  +0.24%  '
Set-TextValue 'B43' 'VeChain'
Set-TextValue 'C43' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D43' '0.0219'
Set-TextValue 'E43' '  +3.84%  '
Set-TextValue 'B44' 'TrustWalletToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D44' '1.20'
Set-TextValue 'E44' '  +2.68%  '
Set-TextValue 'D45' '16.41'
Set-TextValue 'E45' '  +2.18%  '
Set-TextValue 'D46' '1.409.57'
Set-TextValue 'E46' '  +3.01%  '
Set-TextValue 'D47' '91.87'
Set-TextValue 'E47' '  +2.61%  '
Set-TextValue 'E48' '  +2.09%  '
Set-TextValue 'D49' '7.47'
Set-TextValue 'E49' '  +3.47%  '
Set-TextValue 'B50' 'MXToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D50' '2.89'
Set-TextValue 'E50' '  +2.17%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '2.04'
Set-TextValue 'E51' '  +7.21%  '
